# Applies the "what-is-iter8" deck update:
#   - refresh the baked-in "Date Placeholder" field text (6/21/21 -> 10/14/21)
#     on the slide master, every slide layout, and the notes master
#   - reflow / retitle a few shapes on slide 1 (funnel diagram tweaks)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Helper: convert a target EMU value to the Points value that this COM
# host will floor back to the same EMU (1 pt = 12700 EMU). Aiming at the
# middle of the EMU's bucket keeps us safe from float round-trip noise.
# ---------------------------------------------------------------------
function EmuToPt($emu) {
    return ($emu + 0.5) / 12700.0
}

# ---------------------------------------------------------------------
# 1) Update every "Date Placeholder" shape's cached text: 6/21/21 -> 10/14/21
#    (slide master, all slide layouts, notes master)
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
                if ($sh.TextFrame.TextRange.Text -eq "6/21/21") {
                    $sh.TextFrame.TextRange.Text = "10/14/21"
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}

# NOTE: editing NotesMaster shapes directly (Shapes.Item(n).TextFrame...)
# misdirects the write into the slide master in this host, so the notes
# master's date field is updated through its HeadersFooters object instead.
$notesMaster = $p.NotesMaster
$notesMaster.HeadersFooters.DateAndTime.Text = "10/14/21"

# ---------------------------------------------------------------------
# 2) Slide 1 shape tweaks
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)

# "Rectangle 4" - big white backing rectangle grows taller
$rect4 = $s.Shapes.Item("Rectangle 4")
$rect4.Top = EmuToPt 809297
$rect4.Height = EmuToPt 4225157

# "Rectangle 27" - "Assess versions" box widens (shifts left) and its
# text gains a trailing sentence
$rect27 = $s.Shapes.Item("Rectangle 27")
$rect27.Left = EmuToPt 7472855
$rect27.Width = EmuToPt 2002334
$rect27.TextFrame.TextRange.Text = "Assess versions. Find winning version."

# "Rectangle 40" - "Promote winning version" -> "Promote winner"
$rect40 = $s.Shapes.Item("Rectangle 40")
$rect40.TextFrame.TextRange.Text = "Promote winner"

# "Elbow Connector 28" - bent connector shortens to match rect27's new width
$conn28 = $s.Shapes.Item("Elbow Connector 28")
$conn28.Width = EmuToPt 1953673

# "Rectangle 35" - "Deploy new version of application" -> "Deploy new app version"
$rect35 = $s.Shapes.Item("Rectangle 35")
$rect35.TextFrame.TextRange.Text = "Deploy new app version"

# "Elbow Connector 36" - bent connector repositions/resizes
$conn36 = $s.Shapes.Item("Elbow Connector 36")
$conn36.Left = EmuToPt 7888371
$conn36.Top = EmuToPt 2009771
$conn36.Height = EmuToPt 509546
